$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the first
#    (Heading1) paragraph "Play Buffalo Hunter Free - A High
#    Volatility Slot Game".
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$metaRange = $metaPara.Range
$metaRange.End = $metaRange.End - 1
$metaStart = $metaRange.Start
$metaRange.Text = "Meta description: Experience the thrill of Buffalo Hunter with a high jackpot of 12,647 times the bet and variation in free spins. Play for free on mobile or desktop."

$boldLabelRange = $d.Range($metaStart, $metaStart + 16)
$boldLabelRange.Bold = 1

# ------------------------------------------------------------------
# 2) Remove the bold "Play Buffalo Hunter Free - A High Volatility
#    Slot Game" paragraph that used to sit just before the final
#    italic paragraph near the end of the document.
# ------------------------------------------------------------------
$n = $d.Paragraphs.Count
$targetPara = $d.Paragraphs($n - 1)
$targetPara.Range.Delete()

# ------------------------------------------------------------------
# 3) Replace the text of the remaining (italic) final paragraph with
#    the new image-prompt copy, keeping its italic formatting intact.
# ------------------------------------------------------------------
$n2 = $d.Paragraphs.Count
$finalPara = $d.Paragraphs($n2)
$finalRange = $finalPara.Range
$finalRange.End = $finalRange.End - 1
$finalRange.Text = "Create a feature image fitting the game `"Buffalo Hunter`" that features a happy Maya warrior with glasses in a cartoon style. The image should incorporate the North American prairie landscape with wild animals, such as buffaloes, wolves, cougars, and eagles. The warrior should be standing in front of a dream catcher with the game's logo prominently displayed. The overall design should be eye-catching and give players a sense of the game's theme and features."
